$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 566, shifting existing rows 566-642 down to 567-643.
$ws.Range("A566").EntireRow.Insert()

# Populate the newly inserted row with the new weekly price observation.
$ws.Range("A566").Value = 3
$ws.Range("B566").Value = 'Femacal de La Calera'
$ws.Range("C566").Value = 'Coquimbo'
$ws.Range("D566").Value = 45077
$ws.Range("E566").Value = 5
$ws.Range("F566").Value = 'Fruta'
$ws.Range("G566").Value = 100108
$ws.Range("H566").Value = 'Tropicales y subtropicales'
$ws.Range("I566").Value = 100108002
$ws.Range("J566").Value = 'Mango'
$ws.Range("K566").Value = 'Sin especificar'
$ws.Range("L566").Value = 'Primera'
$ws.Range("M566").Value = 228
$ws.Range("N566").Value = 8000
$ws.Range("O566").Value = 8000
$ws.Range("P566").Value = 8000
$ws.Range("Q566").Value = '$/bandeja 4 kilos'
$ws.Range("R566").Value = 'Perú'
$ws.Range("S566").Value = 2000
$ws.Range("T566").Value = 4
